# Update column G ("K") values for rows 2-30 as per the regenerated save_data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    2  = 2
    3  = 0
    4  = 0
    5  = 1
    6  = 1
    7  = 0
    8  = 0
    9  = 2
    10 = 0
    11 = 2
    12 = 0
    13 = 0
    14 = 1
    15 = 0
    16 = 0
    17 = 0
    18 = 1
    19 = 0
    20 = 0
    21 = 1
    22 = 1
    23 = 0
    24 = 0
    25 = 0
    26 = 0
    27 = 2
    28 = 1
    29 = 2
    30 = 1
}

foreach ($row in $newValues.Keys) {
    $ws.Range("G$row").Value = $newValues[$row]
}
